$d = $word.ActiveDocument

# --- Title change ---
$d.Content.Find.Execute("K6 Stress Test Report", $true, $false, $false, $false, $false, $true, 1, $false, "Da3em Stress Test Report", 2) | Out-Null

# --- Update the stress-test results table ---
$table = $d.Tables.Item(1)

# Row 2: 20 VUs
$table.Cell(2, 2).Range.Text = "104.070363"
$table.Cell(2, 3).Range.Text = "1090"

# Row 3: 100 VUs
$table.Cell(3, 2).Range.Text = "290.936517"
$table.Cell(3, 3).Range.Text = "3691"

# Row 4: 500 VUs
$table.Cell(4, 2).Range.Text = "930.300207"
$table.Cell(4, 3).Range.Text = "4754"
$table.Cell(4, 4).Range.Text = "0"

# Row 5: 1000 VUs
$table.Cell(5, 2).Range.Text = "1554.131479"
$table.Cell(5, 3).Range.Text = "5394"
$table.Cell(5, 4).Range.Text = "0"

# New row: 5000 VUs
$newRow = $table.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "5000 VUs"
$newRow.Cells.Item(2).Range.Text = "5364.517892"
$newRow.Cells.Item(3).Range.Text = "9443"
$newRow.Cells.Item(4).Range.Text = "1775"

# --- Update the closing paragraph text ---
$d.Content.Find.Execute("This report includes the essential metrics from the K6 stress test for each stage of the test.", $true, $false, $false, $false, $false, $true, 1, $false, "This report includes the essential metrics from the Da3em K6 stress test for each stage of the test.", 2) | Out-Null
